$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q1" sheet by copying the "2021-Q4" template ---
# (copying preserves header/row styles: bold+border header row, centered index column A)
$template = $wb.Worksheets.Item("2021-Q4")
$beforeTarget = $wb.Worksheets.Item("总计")
$template.Copy($beforeTarget)

# NOTE: sheet references captured before a sheet-collection structural change (add/copy/move)
# can end up pointing at the wrong tab afterwards, since the collection re-indexes. Re-fetch
# both sheets fresh, by name, right after the copy so later writes land on the right tab.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$newSheet = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")

# template has 12 data rows (A2:H12); target needs 17 rows (A2:H17) -> extend 5 more rows
# by cloning the formatting of the last template row downward
for ($r = 13; $r -le 17; $r++) {
    $src = $newSheet.Range("A12:H12")
    $dst = $newSheet.Range("A" + $r + ":H" + $r)
    $src.Copy($dst)
}

# --- Step 2: overwrite header row with the fund-holdings column titles ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Step 3: write the 16 fund rows (A2:H17) ---
# Columns D,E,F,G are numeric-looking percentages/amounts that must stay TEXT (a leading
# apostrophe forces Excel to store them as text, matching the source workbook's t="inlineStr" cells)

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'000727"
$newSheet.Range("C2").Value = "融通健康产业灵活配置混合A"
$newSheet.Range("D2").Value = "'15.30"
$newSheet.Range("E2").Value = "'94.68"
$newSheet.Range("F2").Value = "'9.69"
$newSheet.Range("G2").Value = "'1.4826"
$newSheet.Range("H2").Value = 1

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'002919"
$newSheet.Range("C3").Value = "东吴智慧医疗量化策略灵活配置混合"
$newSheet.Range("D3").Value = "'6.38"
$newSheet.Range("E3").Value = "'93.35"
$newSheet.Range("F3").Value = "'5.36"
$newSheet.Range("G3").Value = "'0.3420"
$newSheet.Range("H3").Value = 8

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'009274"
$newSheet.Range("C4").Value = "融通健康产业灵活配置混合C"
$newSheet.Range("D4").Value = "'3.16"
$newSheet.Range("E4").Value = "'94.68"
$newSheet.Range("F4").Value = "'9.69"
$newSheet.Range("G4").Value = "'0.3062"
$newSheet.Range("H4").Value = 1

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'519087"
$newSheet.Range("C5").Value = "新华优选分红混合"
$newSheet.Range("D5").Value = "'9.03"
$newSheet.Range("E5").Value = "'88.56"
$newSheet.Range("F5").Value = "'3.27"
$newSheet.Range("G5").Value = "'0.2953"
$newSheet.Range("H5").Value = 9

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'000878"
$newSheet.Range("C6").Value = "中海医药健康产业精选灵活配置混合 - A"
$newSheet.Range("D6").Value = "'6.16"
$newSheet.Range("E6").Value = "'89.54"
$newSheet.Range("F6").Value = "'3.46"
$newSheet.Range("G6").Value = "'0.2131"
$newSheet.Range("H6").Value = 10

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'519156"
$newSheet.Range("C7").Value = "新华行业轮换灵活配置混合A"
$newSheet.Range("D7").Value = "'4.98"
$newSheet.Range("E7").Value = "'93.77"
$newSheet.Range("F7").Value = "'3.12"
$newSheet.Range("G7").Value = "'0.1554"
$newSheet.Range("H7").Value = 10

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'519673"
$newSheet.Range("C8").Value = "银河康乐股票"
$newSheet.Range("D8").Value = "'2.31"
$newSheet.Range("E8").Value = "'92.35"
$newSheet.Range("F8").Value = "'4.35"
$newSheet.Range("G8").Value = "'0.1005"
$newSheet.Range("H8").Value = 8

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'006981"
$newSheet.Range("C9").Value = "中金新医药股票A"
$newSheet.Range("D9").Value = "'1.86"
$newSheet.Range("E9").Value = "'91.77"
$newSheet.Range("F9").Value = "'4.74"
$newSheet.Range("G9").Value = "'0.0882"
$newSheet.Range("H9").Value = 3

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'000879"
$newSheet.Range("C10").Value = "中海医药健康产业精选灵活配置混合 - C"
$newSheet.Range("D10").Value = "'2.11"
$newSheet.Range("E10").Value = "'89.54"
$newSheet.Range("F10").Value = "'3.46"
$newSheet.Range("G10").Value = "'0.0730"
$newSheet.Range("H10").Value = 10

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "'001294"
$newSheet.Range("C11").Value = "新华战略新兴产业灵活配置混合"
$newSheet.Range("D11").Value = "'1.07"
$newSheet.Range("E11").Value = "'93.41"
$newSheet.Range("F11").Value = "'4.98"
$newSheet.Range("G11").Value = "'0.0533"
$newSheet.Range("H11").Value = 9

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "'011457"
$newSheet.Range("C12").Value = "新华行业龙头主题股票"
$newSheet.Range("D12").Value = "'0.85"
$newSheet.Range("E12").Value = "'93.55"
$newSheet.Range("F12").Value = "'3.62"
$newSheet.Range("G12").Value = "'0.0308"
$newSheet.Range("H12").Value = 7

$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "'001861"
$newSheet.Range("C13").Value = "富安达健康人生灵活配置混合"
$newSheet.Range("D13").Value = "'0.61"
$newSheet.Range("E13").Value = "'82.18"
$newSheet.Range("F13").Value = "'3.15"
$newSheet.Range("G13").Value = "'0.0192"
$newSheet.Range("H13").Value = 8

$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "'007005"
$newSheet.Range("C14").Value = "中金新医药股票C"
$newSheet.Range("D14").Value = "'0.34"
$newSheet.Range("E14").Value = "'91.77"
$newSheet.Range("F14").Value = "'4.74"
$newSheet.Range("G14").Value = "'0.0161"
$newSheet.Range("H14").Value = 3

$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "'005120"
$newSheet.Range("C15").Value = "上投摩根量化多因子灵活配置混合"
$newSheet.Range("D15").Value = "'0.21"
$newSheet.Range("E15").Value = "'92.54"
$newSheet.Range("F15").Value = "'1.97"
$newSheet.Range("G15").Value = "'0.0041"
$newSheet.Range("H15").Value = 10

$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "'005281"
$newSheet.Range("C16").Value = "中科沃土转型升级灵活配置混合"
$newSheet.Range("D16").Value = "'0.10"
$newSheet.Range("E16").Value = "'21.75"
$newSheet.Range("F16").Value = "'3.49"
$newSheet.Range("G16").Value = "'0.0035"
$newSheet.Range("H16").Value = 2

$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "'519157"
$newSheet.Range("C17").Value = "新华行业轮换灵活配置混合C"
$newSheet.Range("D17").Value = "'0.04"
$newSheet.Range("E17").Value = "'93.77"
$newSheet.Range("F17").Value = "'3.12"
$newSheet.Range("G17").Value = "'0.0012"
$newSheet.Range("H17").Value = 10

# --- Step 4: update the "总计" (totals) sheet - prepend a 2022-Q1 summary row ---
# shift existing data rows 2-6 down to 3-7 (copies cell formatting along with values)
for ($r = 6; $r -ge 2; $r--) {
    $src = $totalSheet.Range("A" + $r + ":D" + $r)
    $dst = $totalSheet.Range("A" + ($r+1) + ":D" + ($r+1))
    $src.Copy($dst)
}

# re-number the shifted rows index column (A) 1..5
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# write the new 2022-Q1 summary row
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 16
$totalSheet.Range("D2").Value = 3.18

